$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I header (row 4): reuse the existing bold/right "year header" style (same as H4) ---
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = 2020

# --- Row 5: new style (fontId=6, General number format, right-aligned, vertical center) ---
$ws.Range("H5").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("I5").NumberFormat = "General"
$ws.Range("I5").HorizontalAlignment = -4152
$ws.Range("I5").Value = 1287.4000000000001

# --- Row 6: new style (fontId=1, General number format, vertical center, no border) ---
$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I6").NumberFormat = "General"
$ws.Range("I6").Value = 56.6

# --- Row 7: same style as row 6/8/9, but left empty ---
$ws.Range("H6").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").NumberFormat = "General"
$ws.Range("I7").ClearContents()

# --- Row 8 ---
$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").NumberFormat = "General"
$ws.Range("I8").Value = 2.5

# --- Row 9 ---
$ws.Range("H9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("I9").NumberFormat = "General"
$ws.Range("I9").Value = 9.3000000000000007

# --- Row 10: new style (fontId=1, General number format, bottom border, vertical center) ---
$ws.Range("H10").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").NumberFormat = "General"
$ws.Range("I10").Value = 0.9

# --- Selection shown in the saved view ---
$ws.Range("L9").Select()
